# Upload Leave Card 12/27/2023 4:01 PM
#
# Applies the author's edits to the "LUNA, GUILLERMA" leave-card workbook:
#   - "2018 LEAVE CREDITS": the monthly accrual schedule (rows 77-119, column A)
#     moves from "1st of month" to "last day of month"; two extra SP (special
#     leave) entries are recorded (rows 85/86 EARNED + row 87 PARTICULARS/REMARKS).
#   - "2017 LEAVE BALANCE": a new VL entry is recorded on row 26.
#   - Active sheet/selection bookkeeping to mirror where the author left off.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("2018 LEAVE CREDITS")
$ws3 = $wb.Worksheets.Item("2017 LEAVE BALANCE")

# ---------------------------------------------------------------------------
# "2018 LEAVE CREDITS" - shift the accrual-period end dates (A77:A119) from
# the 1st of each month to the last day of that same month.
# ---------------------------------------------------------------------------
$periodEnds = [ordered]@{
    77  = 44957
    78  = 44985
    79  = 45016
    80  = 45046
    81  = 45077
    82  = 45107
    83  = 45138
    84  = 45169
    85  = 45199
    86  = 45230
    87  = 45260
    88  = 45291
    89  = 45322
    90  = 45351
    91  = 45382
    92  = 45412
    93  = 45443
    94  = 45473
    95  = 45504
    96  = 45535
    97  = 45565
    98  = 45596
    99  = 45626
    100 = 45657
    101 = 45688
    102 = 45716
    103 = 45747
    104 = 45777
    105 = 45808
    106 = 45838
    107 = 45869
    108 = 45900
    109 = 45930
    110 = 45961
    111 = 45991
    112 = 46022
    113 = 46053
    114 = 46081
    115 = 46112
    116 = 46142
    117 = 46173
    118 = 46203
    119 = 46234
}

foreach ($row in $periodEnds.Keys) {
    $ws2.Cells.Item($row, 1).Value = $periodEnds[$row]
}

# Two new half-day Special Privilege (SP) leave credits taken in Sep/Oct 2023.
$ws2.Range("C85").Value = 1.25
$ws2.Range("C86").Value = 1.25

# A new SP leave entry for Nov 2023 (PARTICULARS + REMARKS / dates taken).
$ws2.Range("B87").Value = "SP(2-0-0)"
$ws2.Range("K87").Value = "11/21,22/2023"

# ---------------------------------------------------------------------------
# "2017 LEAVE BALANCE" - a newly recorded VL availment on row 26.
# ---------------------------------------------------------------------------
$ws3.Range("A26").Value = 45231
$ws3.Range("B26").Value = "VL(3-0-0)"
$ws3.Range("D26").Value = 3

# REMARKS column date needs to render like the other date remarks in this
# column (e.g. K25), not as a plain number.
$ws3.Range("K26").Value = 45086
$ws3.Range("K26").NumberFormat = $ws3.Range("K25").NumberFormat

# ---------------------------------------------------------------------------
# View bookkeeping: the author ended up with "2018 LEAVE CREDITS" active,
# cursor on K87 there; "2017 LEAVE BALANCE" left selected at I9.
# ---------------------------------------------------------------------------
$ws3.Activate()
$ws3.Range("I9").Select()

$ws2.Activate()
$ws2.Range("K87").Select()
